$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 6 (pushes the existing rows 6-34 down to 8-36),
# inheriting formatting (incl. the date style on column D) from row 6.
$ws.Rows.Item(6).Resize(2).Insert()

# New weekly data row for 2023-08-07 (serial 45145), "Primera" quality.
$ws.Cells.Item(6, 1).Value  = 3
$ws.Cells.Item(6, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(6, 3).Value  = "Coquimbo"
$ws.Cells.Item(6, 4).Value  = 45145
$ws.Cells.Item(6, 5).Value  = 5
$ws.Cells.Item(6, 6).Value  = 100112043
$ws.Cells.Item(6, 7).Value  = "Pepino dulce"
$ws.Cells.Item(6, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(6, 9).Value  = "Primera"
$ws.Cells.Item(6, 10).Value = 70
$ws.Cells.Item(6, 11).Value = 22000
$ws.Cells.Item(6, 12).Value = 22000
$ws.Cells.Item(6, 13).Value = 22000
$ws.Cells.Item(6, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(6, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(6, 16).Value = 1467
$ws.Cells.Item(6, 17).Value = 15
$ws.Cells.Item(6, 18).Value = "Hortaliza"

# New weekly data row for 2023-08-07 (serial 45145), "Segunda" quality.
$ws.Cells.Item(7, 1).Value  = 3
$ws.Cells.Item(7, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(7, 3).Value  = "Coquimbo"
$ws.Cells.Item(7, 4).Value  = 45145
$ws.Cells.Item(7, 5).Value  = 5
$ws.Cells.Item(7, 6).Value  = 100112043
$ws.Cells.Item(7, 7).Value  = "Pepino dulce"
$ws.Cells.Item(7, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(7, 9).Value  = "Segunda"
$ws.Cells.Item(7, 10).Value = 60
$ws.Cells.Item(7, 11).Value = 17000
$ws.Cells.Item(7, 12).Value = 17000
$ws.Cells.Item(7, 13).Value = 17000
$ws.Cells.Item(7, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 1133
$ws.Cells.Item(7, 17).Value = 15
$ws.Cells.Item(7, 18).Value = "Hortaliza"
